# Applies odds updates from the 2025-11-25 Betfair Back/Lay daily games sheet update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("AA2").Value = 140
$ws.Range("AE2").Value = 70
$ws.Range("AJ2").Value = 18.5
$ws.Range("AM2").Value = 80
$ws.Range("AO2").Value = 65
$ws.Range("F2").Value = 1.68
$ws.Range("G2").Value = 1.72
$ws.Range("N2").Value = 5.2
$ws.Range("O2").Value = 1.21
$ws.Range("P2").Value = 2.46
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 1.58
$ws.Range("T2").Value = 1.71
$ws.Range("U2").Value = 2.32
$ws.Range("V2").Value = 1.22
$ws.Range("W2").Value = 2.38

# Row 4
$ws.Range("AJ4").Value = 90
$ws.Range("AN4").Value = 42
$ws.Range("H4").Value = 1.92
$ws.Range("I4").Value = 1.95
$ws.Range("P4").Value = 2.24
$ws.Range("S4").Value = 2.92
$ws.Range("V4").Value = 2.04

# Row 5
$ws.Range("AB5").Value = 12
$ws.Range("AC5").Value = 9.800000000000001
$ws.Range("AG5").Value = 9.800000000000001
$ws.Range("AO5").Value = 50
$ws.Range("G5").Value = 1.67
$ws.Range("H5").Value = 5.7
$ws.Range("K5").Value = 4.5
$ws.Range("N5").Value = 5.8
$ws.Range("P5").Value = 2.6
$ws.Range("Q5").Value = 1.58
$ws.Range("T5").Value = 1.66
$ws.Range("W5").Value = 2.48

# Row 6
$ws.Range("F6").Value = 1.45
$ws.Range("I6").Value = 8
$ws.Range("V6").Value = 1.16

# Row 7
$ws.Range("U7").Value = 2.5

# Row 8
$ws.Range("AH8").Value = 28
$ws.Range("AN8").Value = 36
$ws.Range("G8").Value = 2.3
$ws.Range("J8").Value = 3.05
$ws.Range("V8").Value = 1.28
$ws.Range("W8").Value = 1.76
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 34

# Row 9
$ws.Range("Q9").Value = 1.3

# Row 10
$ws.Range("Q10").Value = 1.43

# Row 11
$ws.Range("AB11").Value = 8.4
$ws.Range("AC11").Value = 9.6
$ws.Range("F11").Value = 1.78
$ws.Range("G11").Value = 1.97
$ws.Range("H11").Value = 4.3
$ws.Range("I11").Value = 6.4
$ws.Range("J11").Value = 3.25
$ws.Range("K11").Value = 4.2
$ws.Range("L11").Value = 1.39
$ws.Range("M11").Value = 1.09
$ws.Range("N11").Value = 2.62
$ws.Range("O11").Value = 1.42
$ws.Range("Q11").Value = 2.24
$ws.Range("S11").Value = 3.45
$ws.Range("T11").Value = 2.02
$ws.Range("U11").Value = 1.78
$ws.Range("V11").Value = 1.19
$ws.Range("W11").Value = 2.02

# Row 13
$ws.Range("AC13").Value = 7.6
$ws.Range("AF13").Value = 19.5
$ws.Range("AH13").Value = 15.5
$ws.Range("AI13").Value = 38
$ws.Range("AJ13").Value = 42
$ws.Range("AN13").Value = 24
$ws.Range("U13").Value = 2.34
$ws.Range("X13").Value = 15
$ws.Range("Z13").Value = 17

# Row 14
$ws.Range("AA14").Value = 540
$ws.Range("H14").Value = 13.5
$ws.Range("I14").Value = 14.5
$ws.Range("K14").Value = 7.8
$ws.Range("P14").Value = 3.4
$ws.Range("S14").Value = 1.97
$ws.Range("U14").Value = 2.08

# Row 15
$ws.Range("AA15").Value = 440
$ws.Range("AD15").Value = 40
$ws.Range("AL15").Value = 34
$ws.Range("G15").Value = 1.33
$ws.Range("I15").Value = 12.5
$ws.Range("K15").Value = 6.4
$ws.Range("P15").Value = 2.56
$ws.Range("V15").Value = 1.08

# Row 16
$ws.Range("AF16").Value = 12
$ws.Range("AG16").Value = 9.800000000000001
$ws.Range("AJ16").Value = 20
$ws.Range("I16").Value = 4.6
$ws.Range("J16").Value = 3.9
$ws.Range("N16").Value = 4.6
$ws.Range("S16").Value = 2.94
$ws.Range("X16").Value = 17.5
$ws.Range("Y16").Value = 18.5

# Row 17
$ws.Range("I17").Value = 2.28
$ws.Range("J17").Value = 3.95
$ws.Range("N17").Value = 5.4
$ws.Range("P17").Value = 2.46
$ws.Range("R17").Value = 1.59
$ws.Range("S17").Value = 2.6
$ws.Range("V17").Value = 1.78

# Row 18
$ws.Range("AM18").Value = 40
$ws.Range("AN18").Value = 9.6
$ws.Range("F18").Value = 2.42
$ws.Range("G18").Value = 2.46
$ws.Range("H18").Value = 2.9
$ws.Range("K18").Value = 4.1
$ws.Range("U18").Value = 3.2

# Row 19
$ws.Range("AL19").Value = 55
$ws.Range("AM19").Value = 110
$ws.Range("AN19").Value = 38
$ws.Range("G19").Value = 3.25
$ws.Range("H19").Value = 2.46
$ws.Range("N19").Value = 3.6
$ws.Range("W19").Value = 1.44
